# bom_manual.xlsx - "Updated boms with another alternate part, 0.1uf resistors were OOS"
#
# 1) The 10k resistor group (R3,R8-R19) picked up two more references
#    (R20, R21) now that the 0603 10k part is out of stock and more
#    boards are using the alternate footprint -> qty 13 -> 15.
# 2) A brand new alternate part row is appended for the 0.1uF/C0603
#    capacitor group (same Reference/Value/Footprint as row 4) using a
#    different MPN/vendor part/manufacturer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Row 28: R3.. resistor group gained R20, R21 -----------------------
$ws.Range("A28").Value = 15
$ws.Range("B28").Value = "R3, R8, R9, R10, R11, R12, R13, R14, R15, R16, R17, R18, R19, R20, R21"

# --- 2) New row 54: alternate part for the 0.1u / C_0603 cap group --------
# Clone the formatting of the row directly above (last "Alternates:" entry,
# row 53) so the new row picks up the same fills/number formats, then
# overwrite the values/formulas.
[void]$ws.Range("A53:L53").Copy()
[void]$ws.Range("A54:L54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A54").Value = 27
$ws.Range("B54").Value = "C4, C5, C6, C7, C8, C9, C14, C15, C16, C17, C18, C19, C20, C21, C22, C23, C24, C25, C26, C27, C28, C29, C30, C32, C33, C34, C35"
$ws.Range("C54").Value = "0.1u"
$ws.Range("D54").Value = "Capacitor_SMD:C_0603_1608Metric"
$ws.Range("E54").Value = "FN18F104Z500PSG"
$ws.Range("F54").Value = "JLCPCB "
$ws.Range("G54").Value = "C497032"
$ws.Range("H54").Value = "PSA(Prosperity Dielectrics)"
$ws.Range("I54").Formula = '=IF(F54="JLCPCB ", "", "Yes")'
$ws.Range("J54").Value = 0.003012
$ws.Range("K54").Formula = "=J54*A54"

# Select the freshly added row, matching where the author's cursor ended up
[void]$ws.Rows.Item(54).Select()
